$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3214569873214569
$ws.Range("C2").Value = 4951.46

$ws.Range("A3").Value = 3216549873216549
$ws.Range("C3").Value = 12252.4
